$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object "object[,]" 24,12
$data[0,0] = 0.1412973034146887
$data[0,1] = 0.09374194353000576
$data[0,2] = 0.08542400067607048
$data[0,3] = 1.157840157336764
$data[0,4] = 0.8498590950214577
$data[0,5] = 0.008313253678357307
$data[0,6] = 0.007118001507217642
$data[0,7] = 0.6159309637118611
$data[0,8] = 0.836428130331683
$data[0,9] = 0.1062663075923993
$data[0,10] = 5.034789667622874
$data[0,11] = 0.1007347357716881
$data[1,0] = 0.1236840796445193
$data[1,1] = 0.08242018498263803
$data[1,2] = 0.07637763943655074
$data[1,3] = 1.111567790293336
$data[1,4] = 0.819231082510683
$data[1,5] = 0.01189158735044797
$data[1,6] = 0.01070188559287111
$data[1,7] = 0.6082257122070445
$data[1,8] = 0.8028118093755339
$data[1,9] = 0.09637742208004596
$data[1,10] = 4.387056321386467
$data[1,11] = 0.09159479918989177
$data[2,0] = 0.1127875409493697
$data[2,1] = 0.07550865324168399
$data[2,2] = 0.07081758339515076
$data[2,3] = 1.084023152786514
$data[2,4] = 0.8012721856497791
$data[2,5] = 0.01449370558988577
$data[2,6] = 0.01340849520224108
$data[2,7] = 0.6040138588631976
$data[2,8] = 0.7825683649858775
$data[2,9] = 0.09026090302534229
$data[2,10] = 3.98941619623173
$data[2,11] = 0.08603331138004222
$data[3,0] = 0.1079873676653875
$data[3,1] = 0.07279288764055991
$data[3,2] = 0.06852184997817901
$data[3,3] = 1.0717155328296
$data[3,4] = 0.792933692150072
$data[3,5] = 0.01566013057877846
$data[3,6] = 0.01472571445137083
$data[3,7] = 0.6017961659222664
$data[3,8] = 0.773244954263042
$data[3,9] = 0.0876502860884365
$data[3,10] = 3.827336646592926
$data[3,11] = 0.08388520076105976
$data[4,0] = 0.1067775874288657
$data[4,1] = 0.07245508534058587
$data[4,2] = 0.06810689266587389
$data[4,3] = 1.068112252068801
$data[4,4] = 0.7900730577841557
$data[4,5] = 0.01587025355785487
$data[4,6] = 0.01507221811890513
$data[4,7] = 0.6006716672842458
$data[4,8] = 0.7702786586954176
$data[4,9] = 0.0870869608976399
$data[4,10] = 3.800405981322456
$data[4,11] = 0.0836604047541627
$data[5,0] = 0.1115941253227675
$data[5,1] = 0.07578018663248542
$data[5,2] = 0.07069468939364398
$data[5,3] = 1.079562929182536
$data[5,4] = 0.7970940776800148
$data[5,5] = 0.01453642806055044
$data[5,6] = 0.01373546791342317
$data[5,7] = 0.6018995253488555
$data[5,8] = 0.7785504821053379
$data[5,9] = 0.0898723338944194
$data[5,10] = 3.987183802950085
$data[5,11] = 0.08636344258552953
$data[6,0] = 0.133727189369182
$data[6,1] = 0.09023975860070976
$data[6,2] = 0.08218238710637848
$data[6,3] = 1.13601440754735
$data[6,4] = 0.8337339168376587
$data[6,5] = 0.009488264603959365
$data[6,6] = 0.008606263476932163
$data[6,7] = 0.6104143088418681
$data[6,8] = 0.819590385754303
$data[6,9] = 0.1023936332590409
$data[6,10] = 4.811316508462312
$data[6,11] = 0.0980512199912198
$data[7,0] = 0.1786524571652279
$data[7,1] = 0.1184189505561903
$data[7,2] = 0.1048877226851381
$data[7,3] = 1.262006501744693
$data[7,4] = 0.9203209585986087
$data[7,5] = 0.002984694796962395
$data[7,6] = 0.00238550617196065
$data[7,7] = 0.6352756919829545
$data[7,8] = 0.9106271370464256
$data[7,9] = 0.1272949855884775
$data[7,10] = 6.430359706537388
$data[7,11] = 0.1207920373358959
$data[8,0] = 0.2089560751815611
$data[8,1] = 0.1386022790996009
$data[8,2] = 0.1183852349073753
$data[8,3] = 1.335116206510705
$data[8,4] = 0.9682299309776852
$data[8,5] = 0.00087111863757805
$data[8,6] = 0.0009780706270250406
$data[8,7] = 0.6459614441320696
$data[8,8] = 0.9617624063867609
$data[8,9] = 0.1401590182717882
$data[8,10] = 7.62418160073878
$data[8,11] = 0.1327096371712955
$data[9,0] = 0.2038035364361264
$data[9,1] = 0.1395799445134998
$data[9,2] = 0.09809534266501885
$data[9,3] = 1.170081229133729
$data[9,4] = 0.8245184535077641
$data[9,5] = 0.01928748349961396
$data[9,6] = 0.001430587172760234
$data[9,7] = 0.567134991846558
$data[9,8] = 0.8407692519739314
$data[9,9] = 0.1058574765132896
$data[9,10] = 8.161959027725743
$data[9,11] = 0.09696790796505894
$data[10,0] = 0.1949753257650144
$data[10,1] = 0.1354862088856095
$data[10,2] = 0.08053334847189575
$data[10,3] = 1.029983344415086
$data[10,4] = 0.7076099786486907
$data[10,5] = 0.05800777151707592
$data[10,6] = 0.001429932373506304
$data[10,7] = 0.5053578231777891
$data[10,8] = 0.7400883837264303
$data[10,9] = 0.08171147441368376
$data[10,10] = 8.363090643481314
$data[10,11] = 0.07032964904278671
$data[11,0] = 0.1810489725483961
$data[11,1] = 0.1277634536924523
$data[11,2] = 0.06401778884535148
$data[11,3] = 0.8943197034316839
$data[11,4] = 0.5980627125086073
$data[11,5] = 0.1139880771574724
$data[11,6] = 0.001401807643256703
$data[11,7] = 0.4497377544878134
$data[11,8] = 0.6431967276265098
$data[11,9] = 0.06351920847917647
$data[11,10] = 8.314348499041614
$data[11,11] = 0.04982171226416909
$data[12,0] = 0.168815559132625
$data[12,1] = 0.1208430456242269
$data[12,2] = 0.0533376428274277
$data[12,3] = 0.8023221233174453
$data[12,4] = 0.5255910790403249
$data[12,5] = 0.1636105897294726
$data[12,6] = 0.001511612879227187
$data[12,7] = 0.4139486284501146
$data[12,8] = 0.5776332537986519
$data[12,9] = 0.05447808980800772
$data[12,10] = 8.167101216927847
$data[12,11] = 0.03927449633032687
$data[13,0] = 0.1643777291914006
$data[13,1] = 0.1185337698108668
$data[13,2] = 0.05067595448989337
$data[13,3] = 0.7788306529916866
$data[13,4] = 0.5075218161865962
$data[13,5] = 0.1762651596866789
$data[13,6] = 0.001686813466966974
$data[13,7] = 0.4055594168363257
$data[13,8] = 0.5607521110613121
$data[13,9] = 0.05263295653722988
$data[13,10] = 8.077981516582668
$data[13,11] = 0.0372013881084623
$data[14,0] = 0.1539141134391144
$data[14,1] = 0.1117664862511987
$data[14,2] = 0.04872692472527085
$data[14,3] = 0.7742881568336912
$data[14,4] = 0.5078970109190806
$data[14,5] = 0.1643852024968169
$data[14,6] = 0.002255504775312112
$data[14,7] = 0.4112980716669341
$data[14,8] = 0.5577363713822265
$data[14,9] = 0.051342753739271
$data[14,10] = 7.571662659775484
$data[14,11] = 0.03698526591294105
$data[15,0] = 0.1522519300672371
$data[15,1] = 0.1102533319852057
$data[15,2] = 0.05259129689122233
$data[15,3] = 0.8185162513421744
$data[15,4] = 0.5456541408808704
$data[15,5] = 0.1271691158964217
$data[15,6] = 0.002713472138516693
$data[15,7] = 0.4342873767092641
$data[15,8] = 0.5892951919480183
$data[15,9] = 0.0543779046991606
$data[15,10] = 7.262699657789994
$data[15,11] = 0.04176811113434731
$data[16,0] = 0.1589120843673868
$data[16,1] = 0.1127890729465406
$data[16,2] = 0.06313442938067126
$data[16,3] = 0.916860780044459
$data[16,4] = 0.6267036967271196
$data[16,5] = 0.07449203200684451
$data[16,6] = 0.002619632392396376
$data[16,7] = 0.4780572452129519
$data[16,8] = 0.6600371822846967
$data[16,9] = 0.06507474485436759
$data[16,10] = 7.087006859890721
$data[16,11] = 0.05458031384293349
$data[17,0] = 0.1700923158393834
$data[17,1] = 0.1189578430847007
$data[17,2] = 0.08004802963512603
$data[17,3] = 1.052926368024487
$data[17,4] = 0.7383751825142326
$data[17,5] = 0.02889548436467493
$data[17,6] = 0.002595310230939063
$data[17,7] = 0.5359063778244177
$data[17,8] = 0.75724498648637
$data[17,9] = 0.08614260583649624
$data[17,10] = 7.03032925265785
$data[17,11] = 0.07826705123345334
$data[18,0] = 0.1972970891130359
$data[18,1] = 0.1342212146359856
$data[18,2] = 0.1144067881048301
$data[18,3] = 1.301245257486784
$data[18,4] = 0.9418620153785326
$data[18,5] = 0.00126751721137186
$data[18,6] = 0.001949607277699883
$data[18,7] = 0.6361367315362969
$data[18,8] = 0.9353424027351238
$data[18,9] = 0.1354599797466207
$data[18,10] = 7.309800249224281
$data[18,11] = 0.1304655396048062
$data[19,0] = 0.2243682973293915
$data[19,1] = 0.1511837879452003
$data[19,2] = 0.1298625550548991
$data[19,3] = 1.397922609757501
$data[19,4] = 1.013301425604013
$data[19,5] = 0.0001084192331841649
$data[19,6] = 0.001715053866728766
$data[19,7] = 0.6619664088158999
$data[19,8] = 1.003904720598179
$data[19,9] = 0.1535356276030271
$data[19,10] = 8.229677745702361
$data[19,11] = 0.1478374976836889
$data[20,0] = 0.2426256809394118
$data[20,1] = 0.1614393878051033
$data[20,2] = 0.1384951896869779
$data[20,3] = 1.456400854125093
$data[20,4] = 1.056860241387113
$data[20,5] = 0.00001717757180008661
$data[20,6] = 0.001952629112435211
$data[20,7] = 0.6773094399916317
$data[20,8] = 1.046099038501112
$data[20,9] = 0.163224179190756
$data[20,10] = 8.832845874636519
$data[20,11] = 0.1560672596645603
$data[21,0] = 0.2342615385305749
$data[21,1] = 0.1555817646312079
$data[21,2] = 0.1339966734162878
$data[21,3] = 1.430116328084708
$data[21,4] = 1.038261150069516
$data[21,5] = 0.00001217867005021134
$data[21,6] = 0.001462792442385563
$data[21,7] = 0.6714701187193555
$data[21,8] = 1.028118623434075
$data[21,9] = 0.1584819525744052
$data[21,10] = 8.510890419066072
$data[21,11] = 0.1512231607179473
$data[22,0] = 0.2005171391576965
$data[22,1] = 0.1341651089787916
$data[22,2] = 0.1168754804325083
$data[22,3] = 1.325604562400727
$data[22,4] = 0.9632720115998552
$data[22,5] = 0.001093591912973535
$data[22,6] = 0.0014371354906908
$data[22,7] = 0.647139516579017
$data[22,8] = 0.9543727994392839
$data[22,9] = 0.1398887182644728
$data[22,10] = 7.295621160010285
$data[22,11] = 0.1337750329123537
$data[23,0] = 0.1644480798800174
$data[23,1] = 0.1113126232399111
$data[23,2] = 0.09856563285309505
$data[23,3] = 1.219244329720865
$data[23,4] = 0.8886512039415635
$data[23,5] = 0.004369223442908354
$data[23,6] = 0.004045713917911087
$data[23,7] = 0.6242914289992854
$data[23,8] = 0.8785023940372341
$data[23,9] = 0.1199302050476305
$data[23,10] = 5.991541350702448
$data[23,11] = 0.1152434463997523

$range = $ws.Range("C2:N25")
$range.Value2 = $data
